$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 53.678983
$ws.Range("H2").Value = 161.036949
$ws.Range("I2").Value = 0.1635408689489082
$ws.Range("J2").Value = 0.1635408689489082
$ws.Range("M2").Value = 168.1098273333333
$ws.Range("N2").Value = 504.329482
$ws.Range("O2").Value = 0.2984182258032519
$ws.Range("P2").Value = 0.298418225803252
$ws.Range("Q2").Value = 9023.964563558935
$ws.Range("R2").Value = 81215.68107203041
$ws.Range("S2").Value = 0.04880357595805533
$ws.Range("T2").Value = 0.04880357595805535
$ws.Range("G3").Value = 53.678983
$ws.Range("H3").Value = 161.036949
$ws.Range("I3").Value = 0.1635408689489082
$ws.Range("J3").Value = 0.1635408689489082
$ws.Range("O3").Value = 0.2893586437755394
$ws.Range("P3").Value = 0.2893586437755394
$ws.Range("Q3").Value = 8750.008953244993
$ws.Range("R3").Value = 78750.08057920494
$ws.Range("S3").Value = 0.04732196404092932
$ws.Range("T3").Value = 0.04732196404092932
$ws.Range("G4").Value = 53.678983
$ws.Range("H4").Value = 161.036949
$ws.Range("I4").Value = 0.1635408689489082
$ws.Range("J4").Value = 0.1635408689489082
$ws.Range("M4").Value = 165.99353
$ws.Range("N4").Value = 497.98059
$ws.Range("O4").Value = 0.294661504941043
$ws.Range("P4").Value = 0.294661504941043
$ws.Range("Q4").Value = 8910.363874979988
$ws.Range("R4").Value = 80193.27487481991
$ws.Range("S4").Value = 0.04818919856385119
$ws.Range("T4").Value = 0.04818919856385119
$ws.Range("G5").Value = 53.678983
$ws.Range("H5").Value = 161.036949
$ws.Range("I5").Value = 0.1635408689489082
$ws.Range("J5").Value = 0.1635408689489082
$ws.Range("M5").Value = 66.22673433333334
$ws.Range("N5").Value = 198.680203
$ws.Range("O5").Value = 0.1175616254801657
$ws.Range("P5").Value = 0.1175616254801657
$ws.Range("Q5").Value = 3554.983746424517
$ws.Range("R5").Value = 31994.85371782065
$ws.Range("S5").Value = 0.0192261303860724
$ws.Range("T5").Value = 0.0192261303860724
$ws.Range("I6").Value = 0.327623464087656
$ws.Range("J6").Value = 0.327623464087656
$ws.Range("M6").Value = 168.1098273333333
$ws.Range("N6").Value = 504.329482
$ws.Range("O6").Value = 0.2984182258032519
$ws.Range("P6").Value = 0.298418225803252
$ws.Range("Q6").Value = 18077.82084758923
$ws.Range("R6").Value = 162700.3876283031
$ws.Range("S6").Value = 0.09776881288455375
$ws.Range("T6").Value = 0.09776881288455376
$ws.Range("I7").Value = 0.327623464087656
$ws.Range("J7").Value = 0.327623464087656
$ws.Range("O7").Value = 0.2893586437755394
$ws.Range("P7").Value = 0.2893586437755394
$ws.Range("S7").Value = 0.09480068123744829
$ws.Range("T7").Value = 0.09480068123744829
$ws.Range("I8").Value = 0.327623464087656
$ws.Range("J8").Value = 0.327623464087656
$ws.Range("M8").Value = 165.99353
$ws.Range("N8").Value = 497.98059
$ws.Range("O8").Value = 0.294661504941043
$ws.Range("P8").Value = 0.294661504941043
$ws.Range("Q8").Value = 17850.24317019163
$ws.Range("R8").Value = 160652.1885317247
$ws.Range("S8").Value = 0.09653802298206647
$ws.Range("T8").Value = 0.09653802298206647
$ws.Range("I9").Value = 0.327623464087656
$ws.Range("J9").Value = 0.327623464087656
$ws.Range("M9").Value = 66.22673433333334
$ws.Range("N9").Value = 198.680203
$ws.Range("O9").Value = 0.1175616254801657
$ws.Range("P9").Value = 0.1175616254801657
$ws.Range("Q9").Value = 7121.743312631999
$ws.Range("R9").Value = 64095.68981368799
$ws.Range("S9").Value = 0.03851594698358752
$ws.Range("T9").Value = 0.03851594698358752
$ws.Range("G10").Value = 15.40846566666667
$ws.Range("H10").Value = 46.225397
$ws.Range("I10").Value = 0.04694414319094096
$ws.Range("J10").Value = 0.04694414319094096
$ws.Range("M10").Value = 168.1098273333333
$ws.Range("N10").Value = 504.329482
$ws.Range("O10").Value = 0.2984182258032519
$ws.Range("P10").Value = 0.298418225803252
$ws.Range("Q10").Value = 2590.314502694928
$ws.Range("R10").Value = 23312.83052425435
$ws.Range("S10").Value = 0.01400898792289441
$ws.Range("T10").Value = 0.01400898792289441
$ws.Range("G11").Value = 15.40846566666667
$ws.Range("H11").Value = 46.225397
$ws.Range("I11").Value = 0.04694414319094096
$ws.Range("J11").Value = 0.04694414319094096
$ws.Range("O11").Value = 0.2893586437755394
$ws.Range("P11").Value = 0.2893586437755394
$ws.Range("Q11").Value = 2511.675985722409
$ws.Range("R11").Value = 22605.08387150168
$ws.Range("S11").Value = 0.0135836936069354
$ws.Range("T11").Value = 0.0135836936069354
$ws.Range("G12").Value = 15.40846566666667
$ws.Range("H12").Value = 46.225397
$ws.Range("I12").Value = 0.04694414319094096
$ws.Range("J12").Value = 0.04694414319094096
$ws.Range("M12").Value = 165.99353
$ws.Range("N12").Value = 497.98059
$ws.Range("O12").Value = 0.294661504941043
$ws.Range("P12").Value = 0.294661504941043
$ws.Range("Q12").Value = 2557.705607893803
$ws.Range("R12").Value = 23019.35047104423
$ws.Range("S12").Value = 0.01383263188081048
$ws.Range("T12").Value = 0.01383263188081048
$ws.Range("G13").Value = 15.40846566666667
$ws.Range("H13").Value = 46.225397
$ws.Range("I13").Value = 0.04694414319094096
$ws.Range("J13").Value = 0.04694414319094096
$ws.Range("M13").Value = 66.22673433333334
$ws.Range("N13").Value = 198.680203
$ws.Range("O13").Value = 0.1175616254801657
$ws.Range("P13").Value = 0.1175616254801657
$ws.Range("Q13").Value = 1020.452362190621
$ws.Range("R13").Value = 9184.071259715591
$ws.Range("S13").Value = 0.00551882978030067
$ws.Range("T13").Value = 0.00551882978030067
$ws.Range("G14").Value = 151.606552
$ws.Range("H14").Value = 454.819656
$ws.Range("I14").Value = 0.4618915237724948
$ws.Range("J14").Value = 0.4618915237724948
$ws.Range("M14").Value = 168.1098273333333
$ws.Range("N14").Value = 504.329482
$ws.Range("O14").Value = 0.2984182258032519
$ws.Range("P14").Value = 0.298418225803252
$ws.Range("Q14").Value = 25486.55127932202
$ws.Range("R14").Value = 229378.9615138982
$ws.Range("S14").Value = 0.1378368490377485
$ws.Range("T14").Value = 0.1378368490377485
$ws.Range("G15").Value = 151.606552
$ws.Range("H15").Value = 454.819656
$ws.Range("I15").Value = 0.4618915237724948
$ws.Range("J15").Value = 0.4618915237724948
$ws.Range("O15").Value = 0.2893586437755394
$ws.Range("P15").Value = 0.2893586437755394
$ws.Range("Q15").Value = 24712.81334392275
$ws.Range("R15").Value = 222415.3200953048
$ws.Range("S15").Value = 0.1336523048902264
$ws.Range("T15").Value = 0.1336523048902264
$ws.Range("G16").Value = 151.606552
$ws.Range("H16").Value = 454.819656
$ws.Range("I16").Value = 0.4618915237724948
$ws.Range("J16").Value = 0.4618915237724948
$ws.Range("M16").Value = 165.99353
$ws.Range("N16").Value = 497.98059
$ws.Range("O16").Value = 0.294661504941043
$ws.Range("P16").Value = 0.294661504941043
$ws.Range("Q16").Value = 25165.70673760856
$ws.Range("R16").Value = 226491.3606384771
$ws.Range("S16").Value = 0.1361016515143148
$ws.Range("T16").Value = 0.1361016515143148
$ws.Range("G17").Value = 151.606552
$ws.Range("H17").Value = 454.819656
$ws.Range("I17").Value = 0.4618915237724948
$ws.Range("J17").Value = 0.4618915237724948
$ws.Range("M17").Value = 66.22673433333334
$ws.Range("N17").Value = 198.680203
$ws.Range("O17").Value = 0.1175616254801657
$ws.Range("P17").Value = 0.1175616254801657
$ws.Range("Q17").Value = 10040.40684249669
$ws.Range("R17").Value = 90363.66158247017
$ws.Range("S17").Value = 0.05430071833020506
$ws.Range("T17").Value = 0.05430071833020506
